$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "18-08-2021"
$ws.Range("B5").Value = 6000
$ws.Range("C5").Value = 12420
$ws.Range("D5").Value = 5970
$ws.Range("E5").Value = 1000
$ws.Range("F5").Value = 4970
$ws.Range("G5").Value = 1.54
